$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the diff. Values are written with a leading
# apostrophe to force text interpretation (preventing Excel from auto-
# converting numeric-looking strings like "332.06" or percentages like
# "1.64%" into numbers), then the style is reset to "Normal" so no extra
# cell formatting/style gets introduced.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "332.06"
Set-TextValue $ws.Range("E2") "1.64%"
# Row 3
Set-TextValue $ws.Range("D3") "45.72"
Set-TextValue $ws.Range("E3") "3.89%"
# Row 4
Set-TextValue $ws.Range("D4") "5.695"
Set-TextValue $ws.Range("E4") "3.58%"
# Row 5
Set-TextValue $ws.Range("D5") "0.08358"
Set-TextValue $ws.Range("E5") "4.11%"
# Row 6
Set-TextValue $ws.Range("D6") "2.037"
Set-TextValue $ws.Range("E6") "2.08%"
# Row 7
Set-TextValue $ws.Range("B7") "GateToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.471"
Set-TextValue $ws.Range("E7") "4.30%"
# Row 8
Set-TextValue $ws.Range("B8") "MXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9823"
Set-TextValue $ws.Range("E8") "3.63%"
# Row 9
Set-TextValue $ws.Range("B9") "BTSEToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.597"
Set-TextValue $ws.Range("E9") "0.94%"
# Row 10
Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1170"
Set-TextValue $ws.Range("E10") "1.89%"
# Row 11
Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1943"
Set-TextValue $ws.Range("E11") "5.77%"
# Row 12
Set-TextValue $ws.Range("B12") "MCDex"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D12") "10.38"
Set-TextValue $ws.Range("E12") "-16.99%"
# Row 13
Set-TextValue $ws.Range("B13") "MandalaExchangeToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D13") "0.1006"
Set-TextValue $ws.Range("E13") "3.57%"
# Row 14
Set-TextValue $ws.Range("B14") "BitrueCoin"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.04663"
Set-TextValue $ws.Range("E14") "1.38%"
# Row 15
Set-TextValue $ws.Range("B15") "BitMartToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.1058"
Set-TextValue $ws.Range("E15") "-0.65%"
# Row 16
Set-TextValue $ws.Range("B16") "BitForexToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001283"
Set-TextValue $ws.Range("E16") "0.96%"
# Row 17
Set-TextValue $ws.Range("B17") "TigerCash"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.006045"
Set-TextValue $ws.Range("E17") "5.80%"
# Row 18
Set-TextValue $ws.Range("B18") "LEO"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.370"
Set-TextValue $ws.Range("E18") "0.13%"
# Row 19
Set-TextValue $ws.Range("E19") "-3.68%"
# Row 20
Set-TextValue $ws.Range("E20") "-0.35%"
# Row 21
Set-TextValue $ws.Range("E21") "2.01%"
# Row 22
Set-TextValue $ws.Range("D22") "0.04205"
Set-TextValue $ws.Range("E22") "3.25%"
# Row 23
Set-TextValue $ws.Range("D23") "0.001305"
Set-TextValue $ws.Range("E23") "4.96%"
# Row 24
Set-TextValue $ws.Range("D24") "0.004586"
Set-TextValue $ws.Range("E24") "6.30%"
# Row 25
Set-TextValue $ws.Range("D25") "0.0001280"
Set-TextValue $ws.Range("E25") "7.69%"
# Row 26
Set-TextValue $ws.Range("D26") "0.0003741"
Set-TextValue $ws.Range("E26") "0.09%"
# Row 38
Set-TextValue $ws.Range("D38") "0.02779"
Set-TextValue $ws.Range("E38") "8.38%"
# Row 39
Set-TextValue $ws.Range("D39") "0.05826"
Set-TextValue $ws.Range("E39") "4.78%"
# Row 40
Set-TextValue $ws.Range("E40") "2.59%"
# Row 41
Set-TextValue $ws.Range("D41") "0.1436"
Set-TextValue $ws.Range("E41") "2.95%"
# Row 42
Set-TextValue $ws.Range("D42") "0.007193"
Set-TextValue $ws.Range("E42") "-5.39%"
# Row 43
Set-TextValue $ws.Range("E43") "-1.95%"
# Row 44
Set-TextValue $ws.Range("D44") "0.008179"
Set-TextValue $ws.Range("E44") "-3.95%"
# Row 45
Set-TextValue $ws.Range("D45") "0.00007215"
Set-TextValue $ws.Range("E45") "1.60%"
# Row 46
Set-TextValue $ws.Range("E46") "0.20%"
# Row 47
Set-TextValue $ws.Range("D47") "0.0005802"
Set-TextValue $ws.Range("E47") "-0.17%"
# Row 48
Set-TextValue $ws.Range("E48") "89.23%"
# Row 49
Set-TextValue $ws.Range("D49") "0.003499"
Set-TextValue $ws.Range("E49") "-0.70%"
# Row 50
Set-TextValue $ws.Range("E50") "0.20%"
# Row 51
Set-TextValue $ws.Range("E51") "0.20%"
